$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header: "Escrow Officer"
$ws.Range("E1").Value = "Escrow Officer"

# Fill in the new Escrow Officer test values for rows 2-21 (test1..test20)
for ($i = 1; $i -le 20; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = "test$i"
}

# Size column E to fit its new contents
$ws.Columns("E").ColumnWidth = 13

# Move / update the active selection to F17
$ws.Range("F17").Select() | Out-Null
